$d = $word.ActiveDocument

# The document renames the two logo pictures that are embedded (as
# inline pictures) in the page headers and footers:
#   - footers: Pearson logo picture -> "image2.png"   (was "image1.png")
#   - headers: BTec logo picture    -> "image1.jpg"   (was "image2.jpg")
# This applies to every header/footer instance (first-page and default).
# The pictures are identified by their (stable) alternative text / description
# rather than by their current Name, since Name isn't reliably readable back
# from a freshly-loaded document.

function Rename-InlineLogo($range, $descrMatch, $newName) {
    if ($range.InlineShapes.Count -gt 0) {
        for ($i = 1; $i -le $range.InlineShapes.Count; $i++) {
            $shp = $range.InlineShapes.Item($i)
            if ($shp.AlternativeText -eq $descrMatch) {
                $shp.Name = $newName
            }
        }
    }
}

foreach ($sec in $d.Sections) {
    for ($hfIndex = 1; $hfIndex -le 3; $hfIndex++) {
        $hdr = $sec.Headers.Item($hfIndex)
        if ($hdr.Exists) {
            Rename-InlineLogo $hdr.Range "BTec_Logo-Orange" "image1.jpg"
        }

        $ftr = $sec.Footers.Item($hfIndex)
        if ($ftr.Exists) {
            Rename-InlineLogo $ftr.Range "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" "image2.png"
        }
    }
}
